$d = $word.ActiveDocument

# --- Change 1 ---
# Merge "La fonction 'vérification des donné" + bookmark + "es du formulaire billet' ..."
# into a single run with text "La fonction 'vérification des données du formulaire billet' ..."
# (this also removes the now-stale _GoBack bookmark at this location).
$r1 = $d.Content
$r1.Find.Execute(
    "vérification des données du formulaire billet",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "vérification des données du formulaire billet", 2) | Out-Null

# --- Change 2 ---
# Remove the stray "CB " after "sélectionner sont moyen de paiement".
$r2 = $d.Content
$r2.Find.Execute(
    "sélectionner sont moyen de paiement CB ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sélectionner sont moyen de paiement ", 2) | Out-Null

# Re-insert the _GoBack bookmark (collapsed / zero-length) right after the
# "Paypal" run, i.e. between it and the following spellEnd proofing mark.
$r3 = $d.Content
$r3.Find.Execute("Paypal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r3) | Out-Null
